$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "Alt3" sheet as a copy of "Alt2", placed right
#    after it (becomes the last tab).
# ------------------------------------------------------------------
$alt2 = $wb.Worksheets.Item("Alt2")
$alt2.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$alt3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$alt3.Name = "Alt3"

# ------------------------------------------------------------------
# 2. Replace the REGEX based formula that used to live in B10 with a
#    TEXTSPLIT/MAP based formula (no REGEX functions) that now lives
#    one row lower, in B11 (spilling into B11:B15).
# ------------------------------------------------------------------
$noRegexFormula = '=MAP(B3:B7,LAMBDA(x,TEXTJOIN(", ",,INDEX(TEXTSPLIT(x,{"{","[","*","("},{"}","]","*",")"},1,,),,2))))'

# Clear out the old spilled array formula (this frees up B11:B15 so a
# new formula can be written into B11).
$alt3.Range("B10").Formula2 = ""

$alt3.Range("B11").Formula2 = $noRegexFormula

# Re-apply the original formatting (font/number format) that B10 used
# to carry onto the new spill range B11:B15.
$alt3.Range("B10").Copy()
$alt3.Range("B11:B15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Add the same no-regex formula in column D (spills D10:D14) and
#    give it left/top aligned formatting; D11 additionally wraps text.
# ------------------------------------------------------------------
$alt3.Range("D10").Formula2 = $noRegexFormula

$alt3.Range("D10").HorizontalAlignment = -4131
$alt3.Range("D10").VerticalAlignment = -4160

$alt3.Range("D10").Copy()
$alt3.Range("D11:D14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$alt3.Range("D11").WrapText = $true

# ------------------------------------------------------------------
# 4. Add a demonstration TEXTSPLIT formula in F10 (spills F10:G11)
#    showing how the text is split before being joined above.
# ------------------------------------------------------------------
$splitFormula = '=TEXTSPLIT(B3,{"{","[","*","("},{"}","]","*",")"},1,,)'
$alt3.Range("F10").Formula2 = $splitFormula

# ------------------------------------------------------------------
# 5. Register the hidden AutoFilter defined name for the new sheet,
#    matching the pattern used by Alt1/Alt2/EDA/Original.
# ------------------------------------------------------------------
$alt3.Names.Add("_xlnm._FilterDatabase", "='Alt3'!`$B`$2:`$C`$13", $false)
$fdName = $alt3.Names.Item("_xlnm._FilterDatabase")
$fdName.Visible = $false

# ------------------------------------------------------------------
# 6. Update view state: Alt1 and Alt2 keep cell B10 selected (but are
#    no longer the active tab); Alt3 becomes the active tab with H8
#    selected.
# ------------------------------------------------------------------
$alt1 = $wb.Worksheets.Item("Alt1")
$alt1.Activate()
$alt1.Range("B10").Select()

$alt2.Activate()
$alt2.Range("B10").Select()

$alt3.Activate()
$alt3.Range("H8").Select()
